$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 112
$ws.Range("H112").Value = 3675.2
$ws.Range("I112").Value = 896.6667
$ws.Range("J112").Value = 4054.0908
$ws.Range("K112").Value = 2690.0001
$ws.Range("L112").Value = 12162.2724
$ws.Range("M112").Value = -1582.0001
$ws.Range("N112").Value = -14378.2724

# Row 125
$ws.Range("H125").Value = 1495
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 9000
$ws.Range("M125").Value = -6540

# Row 132
$ws.Range("H132").Value = 2593.6785
$ws.Range("I132").Value = 2542.074
$ws.Range("J132").Value = 3987
$ws.Range("K132").Value = 7626.222
$ws.Range("L132").Value = 11961
$ws.Range("M132").Value = -5096.222
$ws.Range("N132").Value = -17021

# Row 137
$ws.Range("H137").Value = 933
$ws.Range("I137").Value = 766
$ws.Range("J137").Value = 1100
$ws.Range("K137").Value = 2298
$ws.Range("L137").Value = 3300
$ws.Range("M137").Value = 252
$ws.Range("N137").Value = -8400


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 45
$ws.Range("H45").Value = 2248.963
$ws.Range("I45").Value = 2027.3077
$ws.Range("K45").Value = 2027.3077
$ws.Range("M45").Value = -1650.3077

# Row 61
$ws.Range("H61").Value = 2125
$ws.Range("I61").Value = 1833.3334
$ws.Range("K61").Value = 1833.3334
$ws.Range("M61").Value = -1621.3334

# Row 122
$ws.Range("H122").Value = 2346.7058
$ws.Range("I122").Value = 2212.5483
$ws.Range("J122").Value = 3733
$ws.Range("K122").Value = 6637.644899999999
$ws.Range("L122").Value = 11199
$ws.Range("M122").Value = -4187.644899999999
$ws.Range("N122").Value = -16099

# Row 132
$ws.Range("H132").Value = 1957.6086
$ws.Range("I132").Value = 1153.6333
$ws.Range("J132").Value = 3465.0625
$ws.Range("K132").Value = 3460.8999
$ws.Range("L132").Value = 10395.1875
$ws.Range("M132").Value = -930.8998999999999
$ws.Range("N132").Value = -15455.1875

# Row 136
$ws.Range("H136").Value = 2125
$ws.Range("I136").Value = 1833.3334
$ws.Range("K136").Value = 5500.0002
$ws.Range("M136").Value = -2950.0002


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 82
$ws.Range("H82").Value = 18031.572
$ws.Range("J82").Value = 25054.5
$ws.Range("L82").Value = 25054.5
$ws.Range("N82").Value = -25820.5

# Row 85
$ws.Range("H85").Value = 18031.572
$ws.Range("J85").Value = 25054.5
$ws.Range("L85").Value = 25054.5
$ws.Range("N85").Value = -27706.5

# Row 140
$ws.Range("H140").Value = 75969.89999999999
$ws.Range("J140").Value = 79966.664
$ws.Range("L140").Value = 79966.664
$ws.Range("N140").Value = -90326.664


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 2325.2632
$ws.Range("I31").Value = 1961.84
$ws.Range("J31").Value = 3024.1538
$ws.Range("K31").Value = 1961.84
$ws.Range("L31").Value = 3024.1538
$ws.Range("M31").Value = -1666.84
$ws.Range("N31").Value = -3614.1538

# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()

# Row 34
$ws.Range("H34").Value = 2325.2632
$ws.Range("I34").Value = 1961.84
$ws.Range("J34").Value = 3024.1538
$ws.Range("K34").Value = 1961.84
$ws.Range("L34").Value = 3024.1538
$ws.Range("M34").Value = -1759.84
$ws.Range("N34").Value = -3428.1538

# Row 38
$ws.Range("H38").Value = 1875
$ws.Range("J38").Value = 1875
$ws.Range("L38").Value = 1875
$ws.Range("N38").Value = -2629

# Row 44
$ws.Range("H44").Value = 7233.3335
$ws.Range("J44").Value = 7233.3335
$ws.Range("L44").Value = 7233.3335
$ws.Range("N44").Value = -8117.3335

# Row 46
$ws.Range("H46").Value = 1875
$ws.Range("J46").Value = 1875
$ws.Range("L46").Value = 1875
$ws.Range("N46").Value = -2297

# Row 50
$ws.Range("H50").Value = 9234.666999999999
$ws.Range("J50").Value = 9234.666999999999
$ws.Range("L50").Value = 9234.666999999999
$ws.Range("N50").Value = -10484.667

# Row 122
$ws.Range("H122").Value = 1973.2632
$ws.Range("I122").Value = 1779.6923
$ws.Range("J122").Value = 2392.6667
$ws.Range("K122").Value = 5339.0769
$ws.Range("L122").Value = 7178.000100000001
$ws.Range("M122").Value = -2889.0769
$ws.Range("N122").Value = -12078.0001

# Row 132
$ws.Range("H132").Value = 3439.9443
$ws.Range("I132").Value = 2566.6365
$ws.Range("K132").Value = 7699.9095
$ws.Range("M132").Value = -5169.9095


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 20
$ws.Range("H20").Value = 1998.3636
$ws.Range("J20").Value = 1998.3636
$ws.Range("L20").Value = 1998.3636
$ws.Range("N20").Value = -2488.3636

# Row 70
$ws.Range("H70").Value = 9542085
$ws.Range("I70").Value = 17054476
$ws.Range("J70").Value = 7127.654
$ws.Range("K70").Value = 17054476
$ws.Range("L70").Value = 7127.654
$ws.Range("M70").Value = -17054206
$ws.Range("N70").Value = -7667.654

# Row 73
$ws.Range("H73").Value = 9542085
$ws.Range("I73").Value = 17054476
$ws.Range("J73").Value = 7127.654
$ws.Range("K73").Value = 17054476
$ws.Range("L73").Value = 7127.654
$ws.Range("M73").Value = -17053540
$ws.Range("N73").Value = -8999.654

# Row 122
$ws.Range("H122").Value = 2558.6785
$ws.Range("I122").Value = 2562.25
$ws.Range("J122").Value = 2549.75
$ws.Range("K122").Value = 7686.75
$ws.Range("L122").Value = 7649.25
$ws.Range("M122").Value = -5236.75
$ws.Range("N122").Value = -12549.25

# Row 126
$ws.Range("H126").Value = 55556890
$ws.Range("I126").Value = 1995
$ws.Range("J126").Value = 166666670
$ws.Range("K126").Value = 5985
$ws.Range("L126").Value = 500000010
$ws.Range("M126").Value = -3515
$ws.Range("N126").Value = -500004950

# Row 132
$ws.Range("H132").Value = 5027.943
$ws.Range("I132").Value = 5747.44
$ws.Range("J132").Value = 3229.2
$ws.Range("K132").Value = 17242.32
$ws.Range("L132").Value = 9687.599999999999
$ws.Range("M132").Value = -14712.32
$ws.Range("N132").Value = -14747.6

# Row 136
$ws.Range("H136").Value = 10177.75
$ws.Range("J136").Value = 10177.75
$ws.Range("L136").Value = 30533.25
$ws.Range("N136").Value = -35633.25


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 41
$ws.Range("H41").Value = 5099
$ws.Range("J41").Value = 5099
$ws.Range("L41").Value = 5099
$ws.Range("N41").Value = -5975

# Row 136
$ws.Range("H136").Value = 4934.919
$ws.Range("I136").Value = 4456.6523
$ws.Range("J136").Value = 5720.643
$ws.Range("K136").Value = 13369.9569
$ws.Range("L136").Value = 17161.929
$ws.Range("M136").Value = -10819.9569
$ws.Range("N136").Value = -22261.929


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 41
$ws.Range("H41").Value = 21664.857
$ws.Range("J41").Value = 11942.333
$ws.Range("L41").Value = 11942.333
$ws.Range("N41").Value = -12722.333

# Row 136
$ws.Range("H136").Value = 1285.7084
$ws.Range("I136").Value = 1077.8422
$ws.Range("J136").Value = 2075.6
$ws.Range("K136").Value = 3233.5266
$ws.Range("L136").Value = 6226.799999999999
$ws.Range("M136").Value = -683.5266000000001
$ws.Range("N136").Value = -11326.8

